# Auto-generated edit script: updates Leve profit-calculation cells (H:N)
# across multiple worksheets to match the refreshed market-board snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1495   # H40
$ws.Cells.Item(40, 9).Value = 1326.6666   # I40
$ws.Cells.Item(40, 10).Value = 2000   # J40
$ws.Cells.Item(40, 11).Value = 1326.6666   # K40
$ws.Cells.Item(40, 12).Value = 2000   # L40
$ws.Cells.Item(40, 13).Value = -1151.6666   # M40
$ws.Cells.Item(40, 14).Value = -2350   # N40
$ws.Cells.Item(43, 8).Value = 0   # H43
$ws.Cells.Item(43, 9).Value = 0   # I43
$ws.Cells.Item(43, 10).Value = 0   # J43
$ws.Cells.Item(43, 11).Value = 0   # K43
$ws.Cells.Item(43, 12).Value = 0   # L43
$ws.Cells.Item(43, 13).Value = $null   # M43
$ws.Cells.Item(43, 14).Value = $null   # N43
$ws.Cells.Item(98, 8).Value = 9070.200000000001   # H98
$ws.Cells.Item(98, 9).Value = 335   # I98
$ws.Cells.Item(98, 11).Value = 335   # K98
$ws.Cells.Item(98, 13).Value = 1163   # M98
$ws.Cells.Item(100, 8).Value = 997.5   # H100
$ws.Cells.Item(100, 9).Value = 997.5   # I100
$ws.Cells.Item(100, 11).Value = 997.5   # K100
$ws.Cells.Item(100, 13).Value = -456.5   # M100
$ws.Cells.Item(116, 8).Value = 2250   # H116
$ws.Cells.Item(116, 9).Value = 0   # I116
$ws.Cells.Item(116, 11).Value = 0   # K116
$ws.Cells.Item(116, 13).Value = $null   # M116
$ws.Cells.Item(122, 8).Value = 9070.200000000001   # H122
$ws.Cells.Item(122, 9).Value = 335   # I122
$ws.Cells.Item(122, 11).Value = 1005   # K122
$ws.Cells.Item(122, 13).Value = 1445   # M122
$ws.Cells.Item(137, 8).Value = 15407.333   # H137
$ws.Cells.Item(137, 9).Value = 12222   # I137
$ws.Cells.Item(137, 11).Value = 36666   # K137
$ws.Cells.Item(137, 13).Value = -34116   # M137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 723   # H32
$ws.Cells.Item(32, 9).Value = 723   # I32
$ws.Cells.Item(32, 11).Value = 723   # K32
$ws.Cells.Item(32, 13).Value = -436   # M32
$ws.Cells.Item(61, 8).Value = 8377.4   # H61
$ws.Cells.Item(61, 9).Value = 2629   # I61
$ws.Cells.Item(61, 11).Value = 2629   # K61
$ws.Cells.Item(61, 13).Value = -2417   # M61
$ws.Cells.Item(74, 8).Value = 4752.9165   # H74
$ws.Cells.Item(74, 9).Value = 1862.1428   # I74
$ws.Cells.Item(74, 11).Value = 1862.1428   # K74
$ws.Cells.Item(74, 13).Value = -988.1428000000001   # M74
$ws.Cells.Item(77, 8).Value = 4752.9165   # H77
$ws.Cells.Item(77, 9).Value = 1862.1428   # I77
$ws.Cells.Item(77, 11).Value = 9310.714   # K77
$ws.Cells.Item(77, 13).Value = -4942.714   # M77
$ws.Cells.Item(92, 8).Value = 54999   # H92
$ws.Cells.Item(92, 10).Value = 54999   # J92
$ws.Cells.Item(92, 12).Value = 54999   # L92
$ws.Cells.Item(92, 14).Value = -59991   # N92
$ws.Cells.Item(136, 8).Value = 8377.4   # H136
$ws.Cells.Item(136, 9).Value = 2629   # I136
$ws.Cells.Item(136, 11).Value = 7887   # K136
$ws.Cells.Item(136, 13).Value = -5337   # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 49993.6   # H92
$ws.Cells.Item(92, 10).Value = 49993.6   # J92
$ws.Cells.Item(92, 12).Value = 49993.6   # L92
$ws.Cells.Item(92, 14).Value = -54985.6   # N92

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(17, 8).Value = 454   # H17
$ws.Cells.Item(17, 9).Value = 8   # I17
$ws.Cells.Item(17, 10).Value = 900   # J17
$ws.Cells.Item(17, 11).Value = 8   # K17
$ws.Cells.Item(17, 12).Value = 900   # L17
$ws.Cells.Item(17, 13).Value = 166   # M17
$ws.Cells.Item(17, 14).Value = -1248   # N17
$ws.Cells.Item(31, 8).Value = 7999.8   # H31
$ws.Cells.Item(31, 9).Value = 1999.5   # I31
$ws.Cells.Item(31, 11).Value = 1999.5   # K31
$ws.Cells.Item(31, 13).Value = -1704.5   # M31
$ws.Cells.Item(34, 8).Value = 7999.8   # H34
$ws.Cells.Item(34, 9).Value = 1999.5   # I34
$ws.Cells.Item(34, 11).Value = 1999.5   # K34
$ws.Cells.Item(34, 13).Value = -1797.5   # M34
$ws.Cells.Item(58, 8).Value = 7860   # H58
$ws.Cells.Item(58, 10).Value = 11999.667   # J58
$ws.Cells.Item(58, 12).Value = 11999.667   # L58
$ws.Cells.Item(58, 14).Value = -12405.667   # N58
$ws.Cells.Item(86, 8).Value = 10371.75   # H86
$ws.Cells.Item(86, 9).Value = 9662.666999999999   # I86
$ws.Cells.Item(86, 11).Value = 9662.666999999999   # K86
$ws.Cells.Item(86, 13).Value = -8539.666999999999   # M86
$ws.Cells.Item(89, 8).Value = 10371.75   # H89
$ws.Cells.Item(89, 9).Value = 9662.666999999999   # I89
$ws.Cells.Item(89, 11).Value = 48313.335   # K89
$ws.Cells.Item(89, 13).Value = -42697.335   # M89
$ws.Cells.Item(92, 8).Value = 100000   # H92
$ws.Cells.Item(92, 10).Value = 100000   # J92
$ws.Cells.Item(92, 12).Value = 100000   # L92
$ws.Cells.Item(92, 14).Value = -104992   # N92
$ws.Cells.Item(107, 8).Value = 971.2857   # H107
$ws.Cells.Item(107, 9).Value = 999.8889   # I107
$ws.Cells.Item(107, 11).Value = 999.8889   # K107
$ws.Cells.Item(107, 13).Value = 920.1111   # M107
$ws.Cells.Item(132, 8).Value = 12249.5   # H132
$ws.Cells.Item(132, 10).Value = 14999.333   # J132
$ws.Cells.Item(132, 12).Value = 44997.999   # L132
$ws.Cells.Item(132, 14).Value = -50057.999   # N132
$ws.Cells.Item(134, 8).Value = 1999.3334   # H134
$ws.Cells.Item(134, 10).Value = 1999   # J134
$ws.Cells.Item(134, 12).Value = 5997   # L134
$ws.Cells.Item(134, 14).Value = -11067   # N134
$ws.Cells.Item(136, 8).Value = 7860   # H136
$ws.Cells.Item(136, 10).Value = 11999.667   # J136
$ws.Cells.Item(136, 12).Value = 35999.001   # L136
$ws.Cells.Item(136, 14).Value = -41099.001   # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 879.3333   # H8
$ws.Cells.Item(8, 9).Value = 879.3333   # I8
$ws.Cells.Item(8, 11).Value = 2637.9999   # K8
$ws.Cells.Item(8, 13).Value = -2498.9999   # M8
$ws.Cells.Item(12, 8).Value = 217.8   # H12
$ws.Cells.Item(12, 10).Value = 242.83333   # J12
$ws.Cells.Item(12, 12).Value = 728.49999   # L12
$ws.Cells.Item(12, 14).Value = -1074.49999   # N12
$ws.Cells.Item(26, 8).Value = 0   # H26
$ws.Cells.Item(26, 9).Value = 0   # I26
$ws.Cells.Item(26, 10).Value = 0   # J26
$ws.Cells.Item(26, 11).Value = 0   # K26
$ws.Cells.Item(26, 12).Value = 0   # L26
$ws.Cells.Item(26, 13).Value = $null   # M26
$ws.Cells.Item(26, 14).Value = $null   # N26
$ws.Cells.Item(60, 8).Value = 0   # H60
$ws.Cells.Item(60, 9).Value = 0   # I60
$ws.Cells.Item(60, 11).Value = 0   # K60
$ws.Cells.Item(60, 13).Value = $null   # M60
$ws.Cells.Item(113, 8).Value = 374   # H113
$ws.Cells.Item(113, 9).Value = 333   # I113
$ws.Cells.Item(113, 11).Value = 999   # K113
$ws.Cells.Item(113, 13).Value = 1171   # M113
$ws.Cells.Item(117, 8).Value = 441.1111   # H117
$ws.Cells.Item(117, 10).Value = 512   # J117
$ws.Cells.Item(117, 12).Value = 1536   # L117
$ws.Cells.Item(117, 14).Value = -8420   # N117
$ws.Cells.Item(121, 8).Value = 782.1429000000001   # H121
$ws.Cells.Item(121, 9).Value = 674.5   # I121
$ws.Cells.Item(121, 10).Value = 925.6667   # J121
$ws.Cells.Item(121, 11).Value = 2023.5   # K121
$ws.Cells.Item(121, 12).Value = 2777.0001   # L121
$ws.Cells.Item(121, 13).Value = -713.5   # M121
$ws.Cells.Item(121, 14).Value = -5397.0001   # N121

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(9, 8).Value = 721.6667   # H9
$ws.Cells.Item(9, 9).Value = 721.6667   # I9
$ws.Cells.Item(9, 11).Value = 721.6667   # K9
$ws.Cells.Item(9, 13).Value = -497.6667   # M9
$ws.Cells.Item(16, 8).Value = 1540   # H16
$ws.Cells.Item(16, 9).Value = 1540   # I16
$ws.Cells.Item(16, 11).Value = 1540   # K16
$ws.Cells.Item(16, 13).Value = -1370   # M16
$ws.Cells.Item(93, 8).Value = 1982.9333   # H93
$ws.Cells.Item(93, 9).Value = 1457.7142   # I93
$ws.Cells.Item(93, 10).Value = 2442.5   # J93
$ws.Cells.Item(93, 11).Value = 1457.7142   # K93
$ws.Cells.Item(93, 12).Value = 2442.5   # L93
$ws.Cells.Item(93, 13).Value = -209.7141999999999   # M93
$ws.Cells.Item(93, 14).Value = -4938.5   # N93
$ws.Cells.Item(100, 8).Value = 2159.8333   # H100
$ws.Cells.Item(100, 9).Value = 1926.25   # I100
$ws.Cells.Item(100, 11).Value = 1926.25   # K100
$ws.Cells.Item(100, 13).Value = -1385.25   # M100
$ws.Cells.Item(122, 8).Value = 4000   # H122
$ws.Cells.Item(122, 10).Value = 4000   # J122
$ws.Cells.Item(122, 12).Value = 12000   # L122
$ws.Cells.Item(122, 14).Value = -16900   # N122

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1500   # H81
$ws.Cells.Item(81, 9).Value = 1500   # I81
$ws.Cells.Item(81, 11).Value = 3000   # K81
$ws.Cells.Item(81, 13).Value = -1939   # M81
$ws.Cells.Item(84, 8).Value = 1500   # H84
$ws.Cells.Item(84, 9).Value = 1500   # I84
$ws.Cells.Item(84, 11).Value = 15000   # K84
$ws.Cells.Item(84, 13).Value = -9696   # M84
$ws.Cells.Item(92, 8).Value = 19666.334   # H92
$ws.Cells.Item(92, 10).Value = 19666.334   # J92
$ws.Cells.Item(92, 12).Value = 19666.334   # L92
$ws.Cells.Item(92, 14).Value = -24658.334   # N92
$ws.Cells.Item(122, 8).Value = 1704.5   # H122
$ws.Cells.Item(122, 9).Value = 1277.2   # I122
$ws.Cells.Item(122, 11).Value = 3831.6   # K122
$ws.Cells.Item(122, 13).Value = -1381.6   # M122
$ws.Cells.Item(126, 8).Value = 3999.75   # H126
$ws.Cells.Item(126, 10).Value = 4499.5   # J126
$ws.Cells.Item(126, 12).Value = 13498.5   # L126
$ws.Cells.Item(126, 14).Value = -18438.5   # N126
$ws.Cells.Item(132, 8).Value = 5399.8887   # H132
$ws.Cells.Item(132, 9).Value = 2620.8   # I132
$ws.Cells.Item(132, 10).Value = 8873.75   # J132
$ws.Cells.Item(132, 11).Value = 7862.400000000001   # K132
$ws.Cells.Item(132, 12).Value = 26621.25   # L132
$ws.Cells.Item(132, 13).Value = -5332.400000000001   # M132
$ws.Cells.Item(132, 14).Value = -31681.25   # N132
$ws.Cells.Item(136, 8).Value = 7127.4   # H136
$ws.Cells.Item(136, 9).Value = 5385.2354   # I136
$ws.Cells.Item(136, 10).Value = 16999.666   # J136
$ws.Cells.Item(136, 11).Value = 16155.7062   # K136
$ws.Cells.Item(136, 12).Value = 50998.99800000001   # L136
$ws.Cells.Item(136, 13).Value = -13605.7062   # M136
$ws.Cells.Item(136, 14).Value = -56098.99800000001   # N136
